# Added error row handling
# Insert a new "Faulty Rows" worksheet ahead of the existing sheets, with
# the same header row (Fiscal Period, Dist Invoice Date, ...) as the
# "Consolidated" sheet, so faulty/erroring rows can be routed there.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no args inserts a new sheet immediately before the
# currently active sheet ("Consolidated" is active/first), so it lands as
# the new first tab - matching the target layout.
$ws = $wb.Worksheets.Add()
$ws.Name = "Faulty Rows"

$headers = @(
    "Fiscal Period",
    "Dist Invoice Date",
    "Dist Org Group Code",
    "Dist Org Group",
    "Distributor Code",
    "Distributor",
    "Group",
    "Customer Code",
    "Customer",
    "Undisputed Amount"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
